# PT_Sudholz.xlsx — "Changed Sudholz2021 xls file to make .15 change in temp"
#
# Subtract 0.15 from every temperature value in column B (rows 2-118),
# rounding the way Excel's own 15-significant-digit storage precision
# would (so the result matches what typing "=B2-0.15" and recalculating
# in real Excel produces, rather than carrying extra IEEE-754 noise in
# the low bits). Also nudges the saved active-cell selection from B4 to
# B6, matching the author's last selected cell when they saved the file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Round-Excel15([double]$x) {
    if ($x -eq 0) { return 0.0 }
    $absInt = [math]::Floor([math]::Abs($x))
    $intDigits = ([string]$absInt).Length
    $decimals = 15 - $intDigits
    if ($decimals -lt 0) { $decimals = 0 }
    if ($decimals -gt 15) { $decimals = 15 }
    return [math]::Round($x, [int]$decimals)
}

$firstRow = 2
$lastRow = 118
$col = 2  # column B

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $current = $cell.Value2
    $cell.Value2 = Round-Excel15($current - 0.15)
}

# Restore the workbook's saved selection to B6 (was B4).
$ws.Range("B6").Select()
